$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 72331.664
$ws.Range("I3").Value = 40000
$ws.Range("J3").Value = 88497.5
$ws.Range("K3").Value = 40000
$ws.Range("L3").Value = 88497.5
$ws.Range("N3").Value = -88725.5
$ws.Range("M3").Value = -39886

$ws.Range("H19").Value = 1182.3
$ws.Range("I19").Value = 1476.1666
$ws.Range("J19").Value = 741.5
$ws.Range("K19").Value = 1476.1666
$ws.Range("L19").Value = 741.5
$ws.Range("M19").Value = -1301.1666

$ws.Range("H62").Value = 7793.625
$ws.Range("I62").Value = 7143.4287
$ws.Range("J62").Value = 12345
$ws.Range("K62").Value = 7143.4287
$ws.Range("L62").Value = 12345
$ws.Range("M62").Value = -6519.4287

$ws.Range("H65").Value = 7793.625
$ws.Range("I65").Value = 7143.4287
$ws.Range("J65").Value = 12345
$ws.Range("K65").Value = 35717.14350000001
$ws.Range("L65").Value = 61725
$ws.Range("M65").Value = -32597.14350000001

$ws.Range("H86").Value = 333448300
$ws.Range("I86").Value = 666668350
$ws.Range("J86").Value = 166838270
$ws.Range("K86").Value = 666668350
$ws.Range("L86").Value = 166838270
$ws.Range("M86").Value = -666667227

$ws.Range("H89").Value = 333448300
$ws.Range("I89").Value = 666668350
$ws.Range("J89").Value = 166838270
$ws.Range("K89").Value = 3333341750
$ws.Range("L89").Value = 834191350
$ws.Range("M89").Value = -3333336134

$ws.Range("H102").Value = 72331.664
$ws.Range("I102").Value = 40000
$ws.Range("J102").Value = 88497.5
$ws.Range("K102").Value = 40000
$ws.Range("L102").Value = 88497.5
$ws.Range("N102").Value = -94987.5
$ws.Range("M102").Value = -36755

$ws.Range("H121").Value = 2479.2
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2479.2
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 7437.599999999999
$ws.Range("N121").Value = -10931.6

$ws.Range("H132").Value = 1774.762
$ws.Range("I132").Value = 1396.5151
$ws.Range("J132").Value = 3161.6667
$ws.Range("K132").Value = 4189.5453
$ws.Range("L132").Value = 9485.000100000001
$ws.Range("M132").Value = -1659.5453

$ws.Range("H137").Value = 66335.08
$ws.Range("I137").Value = 275793.34
$ws.Range("J137").Value = 3497.6
$ws.Range("K137").Value = 827380.02
$ws.Range("L137").Value = 10492.8
$ws.Range("M137").Value = -824830.02
$ws.Range("N137").Value = -15592.8

$ws.Range("H138").Value = 2576.9592
$ws.Range("I138").Value = 1079.2
$ws.Range("J138").Value = 2747.1592
$ws.Range("K138").Value = 3237.6
$ws.Range("L138").Value = 8241.4776
$ws.Range("M138").Value = 1902.4
$ws.Range("N138").Value = -18521.4776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3523.6667
$ws.Range("I45").Value = 2557.4285
$ws.Range("J45").Value = 3861.85
$ws.Range("K45").Value = 2557.4285
$ws.Range("L45").Value = 3861.85
$ws.Range("M45").Value = -2180.4285
$ws.Range("N45").Value = -4615.85

$ws.Range("H61").Value = 4428.579
$ws.Range("I61").Value = 4524.375
$ws.Range("J61").Value = 4358.909
$ws.Range("K61").Value = 4524.375
$ws.Range("L61").Value = 4358.909
$ws.Range("M61").Value = -4312.375
$ws.Range("N61").Value = -4782.909

$ws.Range("H63").Value = 2951
$ws.Range("I63").Value = 2068.6667
$ws.Range("J63").Value = 4009.8
$ws.Range("K63").Value = 2068.6667
$ws.Range("L63").Value = 4009.8
$ws.Range("M63").Value = -1382.6667

$ws.Range("H66").Value = 2951
$ws.Range("I66").Value = 2068.6667
$ws.Range("J66").Value = 4009.8
$ws.Range("K66").Value = 10343.3335
$ws.Range("L66").Value = 20049
$ws.Range("M66").Value = -6911.333500000001

$ws.Range("H74").Value = 2768.85
$ws.Range("I74").Value = 3293.0833
$ws.Range("J74").Value = 1982.5
$ws.Range("K74").Value = 3293.0833
$ws.Range("L74").Value = 1982.5
$ws.Range("M74").Value = -2419.0833

$ws.Range("H77").Value = 2768.85
$ws.Range("I77").Value = 3293.0833
$ws.Range("J77").Value = 1982.5
$ws.Range("K77").Value = 16465.4165
$ws.Range("L77").Value = 9912.5
$ws.Range("M77").Value = -12097.4165

$ws.Range("H129").Value = 100388
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 100388
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 100388
$ws.Range("N129").Value = -110388

$ws.Range("H132").Value = 360421.03
$ws.Range("I132").Value = 387838.28
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 1163514.84
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -1160984.84

$ws.Range("H136").Value = 4428.579
$ws.Range("I136").Value = 4524.375
$ws.Range("J136").Value = 4358.909
$ws.Range("K136").Value = 13573.125
$ws.Range("L136").Value = 13076.727
$ws.Range("M136").Value = -11023.125
$ws.Range("N136").Value = -18176.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 14000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 14000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 14000
$ws.Range("N44").Value = -14994

$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7840

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 119000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 119000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 119000
$ws.Range("N133").Value = -129120

$ws.Range("H134").Value = 6670178.5
$ws.Range("I134").Value = 6670178.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 20010535.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -20008000.5

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 186750
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 186750
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 186750
$ws.Range("N20").Value = -187222

$ws.Range("H22").Value = 530
$ws.Range("I22").Value = 595
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 595
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -245
$ws.Range("N22").Value = -1100

$ws.Range("H30").Value = 186750
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 186750
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 186750
$ws.Range("N30").Value = -186932

$ws.Range("H31").Value = 4593.587
$ws.Range("I31").Value = 2045.3334
$ws.Range("J31").Value = 5826.613
$ws.Range("K31").Value = 2045.3334
$ws.Range("L31").Value = 5826.613
$ws.Range("M31").Value = -1750.3334
$ws.Range("N31").Value = -6416.613

$ws.Range("H34").Value = 4593.587
$ws.Range("I34").Value = 2045.3334
$ws.Range("J34").Value = 5826.613
$ws.Range("K34").Value = 2045.3334
$ws.Range("L34").Value = 5826.613
$ws.Range("M34").Value = -1843.3334
$ws.Range("N34").Value = -6230.613

$ws.Range("H58").Value = 2574.0508
$ws.Range("I58").Value = 2343.6345
$ws.Range("J58").Value = 4285.7144
$ws.Range("K58").Value = 2343.6345
$ws.Range("L58").Value = 4285.7144
$ws.Range("M58").Value = -2140.6345
$ws.Range("N58").Value = -4691.7144

$ws.Range("H87").Value = 74325.664
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 74325.664
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 74325.664
$ws.Range("N87").Value = -76697.664

$ws.Range("H90").Value = 74325.664
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 74325.664
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 222976.992
$ws.Range("N90").Value = -234832.992

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H128").Value = 186750
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 186750
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 186750
$ws.Range("N128").Value = -196710

$ws.Range("H136").Value = 2574.0508
$ws.Range("I136").Value = 2343.6345
$ws.Range("J136").Value = 4285.7144
$ws.Range("K136").Value = 7030.9035
$ws.Range("L136").Value = 12857.1432
$ws.Range("M136").Value = -4480.9035
$ws.Range("N136").Value = -17957.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2235
$ws.Range("I21").Value = 1180
$ws.Range("J21").Value = 2498.75
$ws.Range("K21").Value = 3540
$ws.Range("L21").Value = 7496.25
$ws.Range("M21").Value = -3367
$ws.Range("N21").Value = -7842.25

$ws.Range("H137").Value = 6776.619
$ws.Range("I137").Value = 1829.5
$ws.Range("J137").Value = 13372.777
$ws.Range("K137").Value = 5488.5
$ws.Range("L137").Value = 40118.331
$ws.Range("M137").Value = -388.5
$ws.Range("N137").Value = -50318.331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 20000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -20560

$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -20996

$ws.Range("H102").Value = 1744
$ws.Range("I102").Value = 1628.6666
$ws.Range("J102").Value = 2349.5
$ws.Range("K102").Value = 1628.6666
$ws.Range("L102").Value = 2349.5
$ws.Range("M102").Value = -6.666600000000017

$ws.Range("H122").Value = 1562.3636
$ws.Range("I122").Value = 1479.125
$ws.Range("J122").Value = 1784.3334
$ws.Range("K122").Value = 4437.375
$ws.Range("L122").Value = 5353.0002
$ws.Range("M122").Value = -1987.375
$ws.Range("N122").Value = -10253.0002

$ws.Range("H126").Value = 2585.7727
$ws.Range("I126").Value = 2309.8948
$ws.Range("J126").Value = 4333
$ws.Range("K126").Value = 6929.6844
$ws.Range("L126").Value = 12999
$ws.Range("M126").Value = -4459.6844

$ws.Range("H132").Value = 4204.4736
$ws.Range("I132").Value = 3891.9285
$ws.Range("J132").Value = 5079.6
$ws.Range("K132").Value = 11675.7855
$ws.Range("L132").Value = 15238.8
$ws.Range("M132").Value = -9145.7855
$ws.Range("N132").Value = -20298.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 13232.608
$ws.Range("I122").Value = 13107.071
$ws.Range("J122").Value = 13427.889
$ws.Range("K122").Value = 39321.213
$ws.Range("L122").Value = 40283.667
$ws.Range("M122").Value = -36871.213
$ws.Range("N122").Value = -45183.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 116659.336
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 116659.336
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 116659.336
$ws.Range("N16").Value = -117243.336

$ws.Range("H62").Value = 5984.2856
$ws.Range("I62").Value = 3125
$ws.Range("J62").Value = 7128
$ws.Range("K62").Value = 3125
$ws.Range("L62").Value = 7128
$ws.Range("M62").Value = -2501
$ws.Range("N62").Value = -8376

$ws.Range("H65").Value = 5984.2856
$ws.Range("I65").Value = 3125
$ws.Range("J65").Value = 7128
$ws.Range("K65").Value = 15625
$ws.Range("L65").Value = 35640
$ws.Range("M65").Value = -12505
$ws.Range("N65").Value = -41880

$ws.Range("H126").Value = 2822.2144
$ws.Range("I126").Value = 2269.6924
$ws.Range("J126").Value = 10005
$ws.Range("K126").Value = 6809.0772
$ws.Range("L126").Value = 30015
$ws.Range("M126").Value = -4339.0772
